$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column "Autorenewal status" after "Action" (before "Renewal Date")
$ws.Columns("W:W").Insert()
$ws.Range("W1").Value = "Autorenewal status"

# Insert new columns "Global Sales" and "HVD code" after "recommitment end date"
# (before "Currency")
$ws.Columns("AI:AJ").Insert()
$ws.Range("AI1").Value = "Global Sales"
$ws.Range("AJ1").Value = "HVD code"

# Refresh the AutoFilter so its range covers the new last column (AW)
$ws.AutoFilterMode = $false
$ws.Range("A1:AW1").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the new AutoFilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$AW`$1"
    }
}

# Reflect the selected view from the saved workbook (matches the AK-column
# selection left behind when the new columns were inserted)
$ws.Range("AK1:AK1048576").Select()
